$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("site_metrics")
$ws.Range("O7").Value = 0.013726789686609
$ws.Range("O17").Value = 0.01151727751821055
$ws.Range("AK17").Value = $true
$ws.Range("O18").Value = 0.01137338214971951
$ws.Range("AK20").Value = $true
$ws.Range("O21").Value = 0.0640732400226012
$ws.Range("AK21").Value = $true
$ws.Range("O22").Value = 0.0783711918509432
$ws.Range("O23").Value = 0.07697079838339536
$ws.Range("AK23").Value = $true
$ws.Range("O24").Value = 0.008191891887501267
$ws.Range("O29").Value = 0.009999163755762241
$ws.Range("O30").Value = 0.01719173062571977
$ws.Range("O40").Value = 0.09338542614434432
$ws.Range("AK42").Value = $true
$ws.Range("O56").Value = 0.2937784768966491
$ws.Range("AK58").Value = $true
$ws.Range("O59").Value = 0.1351165642020508
$ws.Range("O61").Value = 0.1628505517829314
$ws.Range("AK62").Value = $true
$ws.Range("O63").Value = 0.06449326967786709
$ws.Range("N65").Value = 12.5448275862069
$ws.Range("O65").Value = 0.01993759407123781
$ws.Range("Q65").Value = 3.379310344827586
$ws.Range("O73").Value = 0.007274959831971676
$ws.Range("AK78").Value = $true
$ws.Range("O79").Value = 0.06793490328793789
$ws.Range("O83").Value = 0.07042655210741544
$ws.Range("O84").Value = 0.08406351067368563
$ws.Range("O86").Value = 0.0518494331441559
$ws.Range("AK88").Value = $true
$ws.Range("AK91").Value = $true
$ws.Range("O94").Value = 0.01592755503527511
$ws.Range("AK96").Value = $true
$ws.Range("AK99").Value = $true
$ws.Range("O108").Value = 1.074511161406133
$ws.Range("AK119").Value = $true
$ws.Range("O125").Value = 0.2620600464844049
$ws.Range("AK126").Value = $true
$ws.Range("N127").Value = 26.17777777777777
$ws.Range("O127").Value = 0.004312671957069734
$ws.Range("AK128").Value = $true
$ws.Range("AK132").Value = $true
$ws.Range("AK133").Value = $true
$ws.Range("AK135").Value = $true
$ws.Range("AK136").Value = $true

$ws = $wb.Worksheets.Item("mk_duration")
$ws.Range("M4").Value = 0.1280779938814067
$ws.Range("N4").Value = 1.521725027336913
$ws.Range("O4").Value = 0.2019704433497537
$ws.Range("P4").Value = 82
$ws.Range("Q4").Value = 2833.333333333333
$ws.Range("R4").Value = 0.5833333333333333
$ws.Range("S4").Value = 3.833333333333334
$ws.Range("M26").Value = 0.8923205967302188
$ws.Range("N26").Value = 0.1353684133387214
$ws.Range("O26").Value = 0.01970443349753695
$ws.Range("P26").Value = 8
$ws.Range("Q26").Value = 2674
$ws.Range("S26").Value = 14
$ws.Range("K29").Value = "increasing"
$ws.Range("L29").Value = $true
$ws.Range("M29").Value = 0.04803815873298189
$ws.Range("N29").Value = 1.977030712267485
$ws.Range("O29").Value = 0.2610837438423645
$ws.Range("P29").Value = 106
$ws.Range("Q29").Value = 2820.666666666667
$ws.Range("R29").Value = 0.470076726342711
$ws.Range("S29").Value = -1.581074168797954
$ws.Range("M40").Value = 0.6765909906464433
$ws.Range("N40").Value = 0.4171195161282941
$ws.Range("O40").Value = 0.05982905982905983
$ws.Range("P40").Value = 21
$ws.Range("Q40").Value = 2299
$ws.Range("R40").Value = 0.1041666666666666
$ws.Range("S40").Value = 9.3125
$ws.Range("M43").Value = 0.01414805675092889
$ws.Range("N43").Value = 2.453482133447592
$ws.Range("O43").Value = 0.3306878306878307
$ws.Range("P43").Value = 125
$ws.Range("Q43").Value = 2554.333333333333
$ws.Range("R43").Value = 0.3125
$ws.Range("S43").Value = 0.90625
$ws.Range("K46").Value = "no trend"
$ws.Range("L46").Value = $false
$ws.Range("M46").Value = 0.9340231183672678
$ws.Range("N46").Value = -0.08278421786538928
$ws.Range("O46").Value = -0.02205882352941177
$ws.Range("P46").Value = -3
$ws.Range("Q46").Value = 583.6666666666666
$ws.Range("R46").Value = 0
$ws.Range("S46").Value = 23.5
$ws.Range("D65").Value = 0.2441369378490839
$ws.Range("E65").Value = -1.164708671158552
$ws.Range("F65").Value = -0.1551724137931035
$ws.Range("G65").Value = -63
$ws.Range("I65").Value = -0.2012138188608777
$ws.Range("J65").Value = 12.48366013071895
$ws.Range("M65").Value = 0.1478846168671544
$ws.Range("N65").Value = -1.447043942046255
$ws.Range("O65").Value = -0.1885057471264368
$ws.Range("P65").Value = -82
$ws.Range("R65").Value = -0.2333333333333334
$ws.Range("S65").Value = 13.05
$ws.Range("K72").Value = "no trend"
$ws.Range("L72").Value = $false
$ws.Range("M72").Value = 0.3068190873296919
$ws.Range("N72").Value = 1.021919320446587
$ws.Range("O72").Value = 0.1601731601731602
$ws.Range("P72").Value = 37
$ws.Range("Q72").Value = 1241
$ws.Range("R72").Value = 0.425
$ws.Range("S72").Value = 1.787500000000001
$ws.Range("K122").Value = "no trend"
$ws.Range("L122").Value = $false
$ws.Range("M122").Value = 0.2669931288036751
$ws.Range("N122").Value = 1.11001368647762
$ws.Range("O122").Value = 0.1699604743083004
$ws.Range("P122").Value = 43
$ws.Range("Q122").Value = 1431.666666666667
$ws.Range("R122").Value = 0.6
$ws.Range("S122").Value = 8.4
$ws.Range("M127").Value = 0.2837568017823187
$ws.Range("N127").Value = 1.071918139870564
$ws.Range("O127").Value = 0.1985294117647059
$ws.Range("P127").Value = 27
$ws.Range("Q127").Value = 588.3333333333334
$ws.Range("R127").Value = 0.9615384615384616
$ws.Range("S127").Value = 12.80769230769231

$ws = $wb.Worksheets.Item("mk_intra_annual")
$ws.Range("M4").Value = 0.744424299369506
$ws.Range("N4").Value = 0.3260000736978335
$ws.Range("O4").Value = 0.04433497536945813
$ws.Range("P4").Value = 18
$ws.Range("Q4").Value = 2719.333333333333
$ws.Range("M26").Value = 0.796501869623873
$ws.Range("N26").Value = -0.2578769564326593
$ws.Range("O26").Value = -0.03448275862068965
$ws.Range("P26").Value = -14
$ws.Range("Q26").Value = 2541.333333333333
$ws.Range("M29").Value = 0.3489299770046341
$ws.Range("N29").Value = -0.9366668039294629
$ws.Range("O29").Value = -0.1231527093596059
$ws.Range("P29").Value = -50
$ws.Range("Q29").Value = 2736.666666666667
$ws.Range("M40").Value = 0.4749360081817646
$ws.Range("N40").Value = -0.7144709581221618
$ws.Range("O40").Value = -0.09686609686609686
$ws.Range("P40").Value = -34
$ws.Range("Q40").Value = 2133.333333333333
$ws.Range("M43").Value = 0.6658494204997139
$ws.Range("N43").Value = 0.4318513981799204
$ws.Range("O43").Value = 0.0582010582010582
$ws.Range("P43").Value = 22
$ws.Range("Q43").Value = 2364.666666666667
$ws.Range("S43").Value = 3
$ws.Range("K46").Value = "no trend"
$ws.Range("L46").Value = $false
$ws.Range("M46").Value = 0.8885976662382817
$ws.Range("N46").Value = -0.1400788814320553
$ws.Range("O46").Value = -0.02941176470588235
$ws.Range("P46").Value = -4
$ws.Range("Q46").Value = 458.6666666666667
$ws.Range("S46").Value = 1
$ws.Range("D65").Value = 0.07954730198972015
$ws.Range("E65").Value = 1.753318789833714
$ws.Range("F65").Value = 0.2216748768472906
$ws.Range("G65").Value = 90
$ws.Range("H65").Value = 2576.666666666667
$ws.Range("M65").Value = 0.1918233548641619
$ws.Range("N65").Value = 1.305204110593214
$ws.Range("O65").Value = 0.1632183908045977
$ws.Range("P65").Value = 71
$ws.Range("Q65").Value = 2876.333333333333
$ws.Range("M72").Value = 0.7297905219885572
$ws.Range("N72").Value = -0.3454041975020423
$ws.Range("O72").Value = -0.05627705627705628
$ws.Range("P72").Value = -13
$ws.Range("Q72").Value = 1207
$ws.Range("S72").Value = 2
$ws.Range("M122").Value = 0.04249511726808941
$ws.Range("N122").Value = 2.028638567532898
$ws.Range("O122").Value = 0.2885375494071146
$ws.Range("P122").Value = 73
$ws.Range("Q122").Value = 1259.666666666667
$ws.Range("R122").Value = 0.05555555555555555
$ws.Range("S122").Value = 1.388888888888889
$ws.Range("M127").Value = 0.1433090593817012
$ws.Range("N127").Value = 1.463578833494457
$ws.Range("O127").Value = 0.2573529411764706
$ws.Range("P127").Value = 35
$ws.Range("Q127").Value = 539.6666666666666
$ws.Range("R127").Value = 0.08333333333333333
$ws.Range("S127").Value = 1.333333333333333
